$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5532.4546
$ws.Range("J17").Value = 5667.095
$ws.Range("L17").Value = 17001.285
$ws.Range("N17").Value = -17337.285
$ws.Range("H33").Value = 883
$ws.Range("I33").Value = 579.1667
$ws.Range("K33").Value = 579.1667
$ws.Range("M33").Value = -350.1667
$ws.Range("H38").Value = 549.5
$ws.Range("I38").Value = 549.5
$ws.Range("K38").Value = 1648.5
$ws.Range("M38").Value = -1276.5
$ws.Range("H50").Value = 312
$ws.Range("J50").Value = 312
$ws.Range("L50").Value = 936
$ws.Range("N50").Value = -1886
$ws.Range("H58").Value = 377.125
$ws.Range("I58").Value = 416.7143
$ws.Range("J58").Value = 100
$ws.Range("K58").Value = 1250.1429
$ws.Range("L58").Value = 300
$ws.Range("M58").Value = -1100.1429
$ws.Range("N58").Value = -600
$ws.Range("H61").Value = 263
$ws.Range("I61").Value = 263
$ws.Range("K61").Value = 789
$ws.Range("M61").Value = -617
$ws.Range("H62").Value = 8577.111000000001
$ws.Range("I62").Value = 7500.5713
$ws.Range("J62").Value = 12345
$ws.Range("K62").Value = 7500.5713
$ws.Range("L62").Value = 12345
$ws.Range("M62").Value = -6876.5713
$ws.Range("N62").Value = -13593
$ws.Range("H65").Value = 8577.111000000001
$ws.Range("I65").Value = 7500.5713
$ws.Range("J65").Value = 12345
$ws.Range("K65").Value = 37502.85649999999
$ws.Range("L65").Value = 61725
$ws.Range("M65").Value = -34382.85649999999
$ws.Range("N65").Value = -67965
$ws.Range("H82").Value = 4800
$ws.Range("I82").Value = 4800
$ws.Range("K82").Value = 14400
$ws.Range("M82").Value = -13994
$ws.Range("H85").Value = 4800
$ws.Range("I85").Value = 4800
$ws.Range("K85").Value = 14400
$ws.Range("M85").Value = -12996
$ws.Range("H97").Value = 1702.0714
$ws.Range("J97").Value = 1693.7273
$ws.Range("L97").Value = 5081.1819
$ws.Range("N97").Value = -6073.1819
$ws.Range("H99").Value = 780.1429000000001
$ws.Range("J99").Value = 499.25
$ws.Range("L99").Value = 1497.75
$ws.Range("N99").Value = -4493.75
$ws.Range("H101").Value = 675.1429000000001
$ws.Range("J101").Value = 1249.5
$ws.Range("L101").Value = 3748.5
$ws.Range("N101").Value = -6992.5
$ws.Range("H104").Value = 1067.7142
$ws.Range("J104").Value = 1930.5
$ws.Range("L104").Value = 5791.5
$ws.Range("N104").Value = -9285.5
$ws.Range("H112").Value = 57894.777
$ws.Range("I112").Value = 2879.6
$ws.Range("J112").Value = 79054.46000000001
$ws.Range("K112").Value = 8638.799999999999
$ws.Range("L112").Value = 237163.38
$ws.Range("M112").Value = -7530.799999999999
$ws.Range("N112").Value = -239379.38
$ws.Range("H118").Value = 992
$ws.Range("I118").Value = 992
$ws.Range("K118").Value = 2976
$ws.Range("M118").Value = -1319
$ws.Range("H129").Value = 1110.75
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").Value = $null
$ws.Range("H135").Value = 1236.6154
$ws.Range("I135").Value = 807.381
$ws.Range("J135").Value = 3039.4
$ws.Range("K135").Value = 7266.429
$ws.Range("L135").Value = 27354.6
$ws.Range("M135").Value = -4731.429
$ws.Range("N135").Value = -32424.6
$ws.Range("H137").Value = 3526230.8
$ws.Range("J137").Value = 6948396.5
$ws.Range("L137").Value = 20845189.5
$ws.Range("N137").Value = -20850289.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null
$ws.Range("H63").Value = 3062.0715
$ws.Range("I63").Value = 2041
$ws.Range("J63").Value = 4900
$ws.Range("K63").Value = 2041
$ws.Range("L63").Value = 4900
$ws.Range("M63").Value = -1355
$ws.Range("N63").Value = -6272
$ws.Range("H66").Value = 3062.0715
$ws.Range("I66").Value = 2041
$ws.Range("J66").Value = 4900
$ws.Range("K66").Value = 10205
$ws.Range("L66").Value = 24500
$ws.Range("M66").Value = -6773
$ws.Range("N66").Value = -31364
$ws.Range("H115").Value = 85000
$ws.Range("J115").Value = 85000
$ws.Range("L115").Value = 85000
$ws.Range("N115").Value = -88134

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1473.3334
$ws.Range("I107").Value = 1157.5
$ws.Range("K107").Value = 1157.5
$ws.Range("M107").Value = 762.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4362.863
$ws.Range("I31").Value = 1995.0588
$ws.Range("K31").Value = 1995.0588
$ws.Range("M31").Value = -1700.0588
$ws.Range("H34").Value = 4362.863
$ws.Range("I34").Value = 1995.0588
$ws.Range("K34").Value = 1995.0588
$ws.Range("M34").Value = -1793.0588
$ws.Range("H132").Value = 3357.6924
$ws.Range("I132").Value = 2387.5
$ws.Range("K132").Value = 7162.5
$ws.Range("M132").Value = -4632.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3533.3333
$ws.Range("I17").Value = 10000
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 900
$ws.Range("M17").Value = -29831
$ws.Range("N17").Value = -1238
$ws.Range("H23").Value = 132.08333
$ws.Range("J23").Value = 124.111115
$ws.Range("L23").Value = 372.333345
$ws.Range("N23").Value = -842.333345
$ws.Range("H132").Value = 1251795.9
$ws.Range("J132").Value = 1667766.6
$ws.Range("L132").Value = 15009899.4
$ws.Range("N132").Value = -15014959.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 65934.664
$ws.Range("I113").Value = 48902.75
$ws.Range("K113").Value = 48902.75
$ws.Range("M113").Value = -46732.75
$ws.Range("H114").Value = 113500
$ws.Range("J114").Value = 113500
$ws.Range("L114").Value = 113500
$ws.Range("N114").Value = -122178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2348.1
$ws.Range("I22").Value = 1181.8334
$ws.Range("J22").Value = 4097.5
$ws.Range("K22").Value = 1181.8334
$ws.Range("L22").Value = 4097.5
$ws.Range("M22").Value = -886.8334
$ws.Range("N22").Value = -4687.5
$ws.Range("H27").Value = 2348.1
$ws.Range("I27").Value = 1181.8334
$ws.Range("J27").Value = 4097.5
$ws.Range("K27").Value = 1181.8334
$ws.Range("L27").Value = 4097.5
$ws.Range("M27").Value = -1074.8334
$ws.Range("N27").Value = -4311.5
$ws.Range("H110").Value = 98274.39999999999
$ws.Range("J110").Value = 98274.39999999999
$ws.Range("L110").Value = 98274.39999999999
$ws.Range("N110").Value = -106454.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 39623.5
$ws.Range("J74").Value = 45664.668
$ws.Range("L74").Value = 45664.668
$ws.Range("N74").Value = -47536.668
$ws.Range("H77").Value = 39623.5
$ws.Range("J77").Value = 45664.668
$ws.Range("L77").Value = 136994.004
$ws.Range("N77").Value = -146354.004
$ws.Range("H96").Value = 9478.583000000001
$ws.Range("I96").Value = 4799.4287
$ws.Range("J96").Value = 16029.4
$ws.Range("K96").Value = 4799.4287
$ws.Range("L96").Value = 16029.4
$ws.Range("M96").Value = -3426.4287
$ws.Range("N96").Value = -18775.4
